$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed date) column C for rows 2-5 from 2023-10-05 to 2023-10-08
$newDate = Get-Date -Year 2023 -Month 10 -Day 8 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

$ws.Range("C2").Value = $newDate
$ws.Range("C3").Value = $newDate
$ws.Range("C4").Value = $newDate
$ws.Range("C5").Value = $newDate
